$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 2877.4055
$ws.Range("I132").Value = 2497.5334
$ws.Range("J132").Value = 4505.4287
$ws.Range("K132").Value = 7492.600199999999
$ws.Range("L132").Value = 13516.2861
$ws.Range("M132").Value = -4962.600199999999
$ws.Range("N132").Value = -18576.2861
# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 2131570.5
$ws.Range("I137").Value = 5556862
$ws.Range("J137").Value = 5527.6553
$ws.Range("K137").Value = 16670586
$ws.Range("L137").Value = 16582.9659
$ws.Range("M137").Value = -16668036
$ws.Range("N137").Value = -21682.9659

$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 32324790
$ws.Range("I61").Value = 47667956
$ws.Range("J61").Value = 104142.8
$ws.Range("K61").Value = 47667956
$ws.Range("L61").Value = 104142.8
$ws.Range("M61").Value = -47667744
$ws.Range("N61").Value = -104566.8
# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 49624
$ws.Range("I132").Value = 31985.484
$ws.Range("J132").Value = 114298.555
$ws.Range("K132").Value = 95956.452
$ws.Range("L132").Value = 342895.665
$ws.Range("M132").Value = -93426.452
$ws.Range("N132").Value = -347955.665
# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 32324790
$ws.Range("I136").Value = 47667956
$ws.Range("J136").Value = 104142.8
$ws.Range("K136").Value = 143003868
$ws.Range("L136").Value = 312428.4
$ws.Range("M136").Value = -143001318
$ws.Range("N136").Value = -317528.4

$ws = $wb.Worksheets.Item("BSM")
# Row 59: Pop That Top / Cobalt Raising Hammer
$ws.Range("H59").Value = 50000
$ws.Range("J59").Value = 50000
$ws.Range("L59").Value = 50000
$ws.Range("N59").Value = -51694
# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 1855.8511
$ws.Range("I134").Value = 1763.5814
$ws.Range("J134").Value = 2847.75
$ws.Range("K134").Value = 5290.7442
$ws.Range("L134").Value = 8543.25
$ws.Range("M134").Value = -2755.7442
$ws.Range("N134").Value = -13613.25

$ws = $wb.Worksheets.Item("CRP")
# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 35715670
$ws.Range("I58").Value = 47620428
$ws.Range("J58").Value = 1400.5714
$ws.Range("K58").Value = 47620428
$ws.Range("L58").Value = 1400.5714
$ws.Range("M58").Value = -47620225
$ws.Range("N58").Value = -1806.5714
# Row 99: O Pine / Pine Lumber
$ws.Range("H99").Value = 1284.3334
$ws.Range("I99").Value = 1073.1428
$ws.Range("J99").Value = 1580
$ws.Range("K99").Value = 1073.1428
$ws.Range("L99").Value = 1580
$ws.Range("M99").Value = 424.8571999999999
$ws.Range("N99").Value = -4576
# Row 122: Timber of Tenkonto / Horse Chestnut Lumber
$ws.Range("H122").Value = 1324.973
$ws.Range("I122").Value = 1290.3928
$ws.Range("J122").Value = 1432.5555
$ws.Range("K122").Value = 3871.1784
$ws.Range("L122").Value = 4297.666499999999
$ws.Range("M122").Value = -1421.1784
$ws.Range("N122").Value = -9197.666499999999
# Row 126: A Better Conductor / Red Pine Lumber
$ws.Range("H126").Value = 1284.3334
$ws.Range("I126").Value = 1073.1428
$ws.Range("J126").Value = 1580
$ws.Range("K126").Value = 3219.4284
$ws.Range("L126").Value = 4740
$ws.Range("M126").Value = -749.4284000000002
$ws.Range("N126").Value = -9680
# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 18995.158
$ws.Range("I132").Value = 1122.875
$ws.Range("K132").Value = 3368.625
$ws.Range("M132").Value = -838.625
# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 35715670
$ws.Range("I136").Value = 47620428
$ws.Range("J136").Value = 1400.5714
$ws.Range("K136").Value = 142861284
$ws.Range("L136").Value = 4201.7142
$ws.Range("M136").Value = -142858734
$ws.Range("N136").Value = -9301.7142

$ws = $wb.Worksheets.Item("CUL")
# Row 26: A Grape Idea / Grape Juice
$ws.Range("H26").Value = 178881.42
$ws.Range("I26").Value = 120
$ws.Range("J26").Value = 357642.84
$ws.Range("K26").Value = 360
$ws.Range("L26").Value = 1072928.52
$ws.Range("M26").Value = -72
$ws.Range("N26").Value = -1073504.52
# Row 92: Oh No Udon / Gyr Abanian Flour
$ws.Range("H92").Value = 669.9
$ws.Range("I92").Value = 660
$ws.Range("J92").Value = 693
$ws.Range("K92").Value = 1980
$ws.Range("L92").Value = 2079
$ws.Range("M92").Value = -732
$ws.Range("N92").Value = -4575
# Row 107: Slippery Service / Frantoio Oil
$ws.Range("H107").Value = 568.40814
$ws.Range("I107").Value = 458.4074
$ws.Range("J107").Value = 703.4091
$ws.Range("K107").Value = 1375.2222
$ws.Range("L107").Value = 2110.2273
$ws.Range("M107").Value = 544.7778000000001
$ws.Range("N107").Value = -5950.2273
# Row 108: Meet for Meat / Grilled Rail
$ws.Range("H108").Value = 799.25
$ws.Range("I108").Value = 342
$ws.Range("J108").Value = 4000
$ws.Range("K108").Value = 1026
$ws.Range("L108").Value = 12000
$ws.Range("M108").Value = 1854
$ws.Range("N108").Value = -17760
# Row 110: His Dark Utensils / Spaghetti al Nero
$ws.Range("H110").Value = 2756.75
$ws.Range("I110").Value = 1009
$ws.Range("J110").Value = 8000
$ws.Range("K110").Value = 3027
$ws.Range("L110").Value = 24000
$ws.Range("M110").Value = 1063
$ws.Range("N110").Value = -32180
# Row 111: Soup for the Soldier / Broad Bean Soup
$ws.Range("H111").Value = 2873.625
$ws.Range("I111").Value = 498.16666
$ws.Range("J111").Value = 10000
$ws.Range("K111").Value = 1494.49998
$ws.Range("L111").Value = 30000
$ws.Range("M111").Value = 1572.50002
$ws.Range("N111").Value = -36134
# Row 116: On a Full Stomach / Sausage Links
$ws.Range("H116").Value = 122170.91
$ws.Range("I116").Value = 133436
$ws.Range("J116").Value = 112783.336
$ws.Range("K116").Value = 400308
$ws.Range("L116").Value = 338350.008
$ws.Range("M116").Value = -396866
$ws.Range("N116").Value = -345234.008
# Row 119: Super Dark Times / Risotto al Nero
$ws.Range("H119").Value = 3431.8333
$ws.Range("I119").Value = 1524.8182
$ws.Range("J119").Value = 6428.5713
$ws.Range("K119").Value = 4574.4546
$ws.Range("L119").Value = 19285.7139
$ws.Range("M119").Value = 263.5454
$ws.Range("N119").Value = -28961.7139
# Row 121: A Cookie for Your Troubles / Coffee Biscuit
$ws.Range("H121").Value = 51225970
$ws.Range("I121").Value = 1270
$ws.Range("J121").Value = 63590550
$ws.Range("K121").Value = 3810
$ws.Range("L121").Value = 190771650
$ws.Range("M121").Value = -2500
$ws.Range("N121").Value = -190774270
# Row 123: Topping Up the Pot / Zurek
$ws.Range("H123").Value = 2119.6365
$ws.Range("I123").Value = 924.2
$ws.Range("J123").Value = 3115.8333
$ws.Range("K123").Value = 2772.6
$ws.Range("L123").Value = 9347.499899999999
$ws.Range("M123").Value = -322.6000000000004
$ws.Range("N123").Value = -14247.4999
# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 983.9464
$ws.Range("J131").Value = 983.9464
$ws.Range("L131").Value = 2951.8392
$ws.Range("N131").Value = -13031.8392

$ws = $wb.Worksheets.Item("GSM")
# Row 102: Put the Metal to the Peddle / Durium Ingot
$ws.Range("H102").Value = 1939.6
$ws.Range("I102").Value = 1862
$ws.Range("J102").Value = 1991.3334
$ws.Range("K102").Value = 1862
$ws.Range("L102").Value = 1991.3334
$ws.Range("M102").Value = -240
$ws.Range("N102").Value = -5235.3334
# Row 113: Copious Crystal Cannons / Manasilver Nugget
$ws.Range("H113").Value = 2002.826
$ws.Range("I113").Value = 1306.5
$ws.Range("J113").Value = 2762.4546
$ws.Range("K113").Value = 1306.5
$ws.Range("L113").Value = 2762.4546
$ws.Range("M113").Value = 863.5
$ws.Range("N113").Value = -7102.4546
# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 60570.59
$ws.Range("I132").Value = 37321.715
$ws.Range("J132").Value = 169065.33
$ws.Range("K132").Value = 111965.145
$ws.Range("L132").Value = 507195.99
$ws.Range("M132").Value = -109435.145
$ws.Range("N132").Value = -512255.99

$ws = $wb.Worksheets.Item("LTW")
# Row 61: Spelling Me Softly / Raptor Leather
$ws.Range("H61").Value = 2662.7368
$ws.Range("I61").Value = 2414.4285
$ws.Range("K61").Value = 2414.4285
$ws.Range("M61").Value = -2212.4285
# Row 82: Trainin' the Neck / Dragon Leather
$ws.Range("H82").Value = 1294.7
$ws.Range("I82").Value = 1007.8333
$ws.Range("K82").Value = 1007.8333
$ws.Range("M82").Value = -646.8333
# Row 85: Training Is Only Skintight (L) / Dragon Leather
$ws.Range("H85").Value = 1294.7
$ws.Range("I85").Value = 1007.8333
$ws.Range("K85").Value = 1007.8333
$ws.Range("M85").Value = 240.1667
# Row 93: Hide to Go Seek / Gagana Leather
$ws.Range("H93").Value = 1450.5
$ws.Range("I93").Value = 1240.6
$ws.Range("K93").Value = 1240.6
$ws.Range("M93").Value = 7.400000000000091
# Row 100: Tiger in the Sack / Tiger Leather
$ws.Range("H100").Value = 1858.8667
$ws.Range("J100").Value = 1972.5
$ws.Range("L100").Value = 1972.5
$ws.Range("N100").Value = -3054.5
# Row 113: Peace in Rest / Atrociraptor Leather
$ws.Range("H113").Value = 2662.7368
$ws.Range("I113").Value = 2414.4285
$ws.Range("K113").Value = 2414.4285
$ws.Range("M113").Value = -244.4285
# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 58909.305
$ws.Range("I136").Value = 40393.703
$ws.Range("J136").Value = 114456.11
$ws.Range("K136").Value = 121181.109
$ws.Range("L136").Value = 343368.33
$ws.Range("M136").Value = -118631.109
$ws.Range("N136").Value = -348468.33

$ws = $wb.Worksheets.Item("WVR")
# Row 81: Where the Dragonflies, the Net Catches / Crawler Silk
$ws.Range("H81").Value = 2188.0557
$ws.Range("I81").Value = 763.3333
$ws.Range("K81").Value = 1526.6666
$ws.Range("M81").Value = -465.6666
# Row 84: To Kill a Dragon on Nameday (L) / Crawler Silk
$ws.Range("H84").Value = 2188.0557
$ws.Range("I84").Value = 763.3333
$ws.Range("K84").Value = 7633.333000000001
$ws.Range("M84").Value = -2329.333000000001
# Row 113: A Tender Table / Pixie Floss
$ws.Range("H113").Value = 879.5
$ws.Range("I113").Value = 905.2143
$ws.Range("J113").Value = 699.5
$ws.Range("K113").Value = 2715.6429
$ws.Range("L113").Value = 2098.5
$ws.Range("M113").Value = -545.6428999999998
$ws.Range("N113").Value = -6438.5
# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 28031.338
$ws.Range("I132").Value = 17943.518
$ws.Range("J132").Value = 64599.688
$ws.Range("K132").Value = 53830.554
$ws.Range("L132").Value = 193799.064
$ws.Range("M132").Value = -51300.554
$ws.Range("N132").Value = -198859.064
# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 38169.184
$ws.Range("I136").Value = 27165.236
$ws.Range("J136").Value = 62766.234
$ws.Range("K136").Value = 81495.708
$ws.Range("L136").Value = 188298.702
$ws.Range("M136").Value = -78945.708
$ws.Range("N136").Value = -193398.702
